$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2021-11-03"
$ws.Range("B1").Value = "November 2021 (through November 03)"

$ws.Range("AI4").Value = 3
$ws.Range("BE5").Value = 2
$ws.Range("M6").Value = 1
$ws.Range("X7").Value = 1
$ws.Range("M8").Value = 1
$ws.Range("AT9").Value = 2
$ws.Range("M12").Value = 1
$ws.Range("AI15").Value = 1
$ws.Range("M21").Value = 1
$ws.Range("B27").Value = 1
$ws.Range("AI28").Value = 3
$ws.Range("B30").Value = 1
$ws.Range("M32").Value = 1
$ws.Range("AI36").Value = 1
$ws.Range("AT36").Value = 1
$ws.Range("B43").Value = 1
$ws.Range("B49").Value = 1
$ws.Range("AI64").Value = 1
$ws.Range("M68").Value = 1
$ws.Range("AT68").Value = 1
$ws.Range("B79").Value = 1
$ws.Range("M88").Value = 1
$ws.Range("M98").Value = 1
